$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.31746133333333
$ws.Range("H2").Value = 126.952384
$ws.Range("I2").Value = 0.6904142182914543
$ws.Range("J2").Value = 0.6904142182914543
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 83.40125033333334
$ws.Range("N2").Value = 250.203751
$ws.Range("O2").Value = 0.9428346765536562
$ws.Range("P2").Value = 0.9428346765536562
$ws.Range("Q2").Value = 3529.329186132487
$ws.Range("R2").Value = 31763.96267519238
$ws.Range("S2").Value = 0.6509464661908687
$ws.Range("T2").Value = 0.6509464661908687

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.31746133333333
$ws.Range("H3").Value = 126.952384
$ws.Range("I3").Value = 0.6904142182914543
$ws.Range("J3").Value = 0.6904142182914543
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.673647
$ws.Range("N3").Value = 8.020941
$ws.Range("O3").Value = 0.03022505171551549
$ws.Range("P3").Value = 0.03022505171551549
$ws.Range("Q3").Value = 113.1419535414827
$ws.Range("R3").Value = 1018.277581873344
$ws.Range("S3").Value = 0.02086780545298641
$ws.Range("T3").Value = 0.02086780545298641

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.31746133333333
$ws.Range("H4").Value = 126.952384
$ws.Range("I4").Value = 0.6904142182914543
$ws.Range("J4").Value = 0.6904142182914543
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2062613333333333
$ws.Range("N4").Value = 0.618784
$ws.Range("O4").Value = 0.002331743669568637
$ws.Range("P4").Value = 0.002331743669568637
$ws.Range("Q4").Value = 8.72845599789511
$ws.Range("R4").Value = 78.556103981056
$ws.Range("S4").Value = 0.001609868982881278
$ws.Range("T4").Value = 0.001609868982881278

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 42.31746133333333
$ws.Range("H5").Value = 126.952384
$ws.Range("I5").Value = 0.6904142182914543
$ws.Range("J5").Value = 0.6904142182914543
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.881585
$ws.Range("N5").Value = 5.644755
$ws.Range("O5").Value = 0.02127094711161878
$ws.Range("P5").Value = 0.02127094711161878
$ws.Range("Q5").Value = 79.62390048288
$ws.Range("R5").Value = 716.6151043459199
$ws.Range("S5").Value = 0.01468576432238715
$ws.Range("T5").Value = 0.01468576432238715

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 42.31746133333333
$ws.Range("H6").Value = 126.952384
$ws.Range("I6").Value = 0.6904142182914543
$ws.Range("J6").Value = 0.6904142182914543
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2952356666666667
$ws.Range("N6").Value = 0.885707
$ws.Range("O6").Value = 0.003337580949640955
$ws.Range("P6").Value = 0.003337580949640955
$ws.Range("Q6").Value = 12.49362390838756
$ws.Range("R6").Value = 112.442615175488
$ws.Range("S6").Value = 0.00230431334233081
$ws.Range("T6").Value = 0.00230431334233081

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.08476666666667
$ws.Range("H7").Value = 33.2543
$ws.Range("I7").Value = 0.1808492350906109
$ws.Range("J7").Value = 0.1808492350906109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 83.40125033333334
$ws.Range("N7").Value = 250.203751
$ws.Range("O7").Value = 0.9428346765536562
$ws.Range("P7").Value = 0.9428346765536562
$ws.Range("Q7").Value = 924.4833996532556
$ws.Range("R7").Value = 8320.3505968793
$ws.Range("S7").Value = 0.1705109300716323
$ws.Range("T7").Value = 0.1705109300716323

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.08476666666667
$ws.Range("H8").Value = 33.2543
$ws.Range("I8").Value = 0.1808492350906109
$ws.Range("J8").Value = 0.1808492350906109
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.673647
$ws.Range("N8").Value = 8.020941
$ws.Range("O8").Value = 0.03022505171551549
$ws.Range("P8").Value = 0.03022505171551549
$ws.Range("Q8").Value = 29.63675314403334
$ws.Range("R8").Value = 266.7307782963
$ws.Range("S8").Value = 0.005466177483325133
$ws.Range("T8").Value = 0.005466177483325133

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.08476666666667
$ws.Range("H9").Value = 33.2543
$ws.Range("I9").Value = 0.1808492350906109
$ws.Range("J9").Value = 0.1808492350906109
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2062613333333333
$ws.Range("N9").Value = 0.618784
$ws.Range("O9").Value = 0.002331743669568637
$ws.Range("P9").Value = 0.002331743669568637
$ws.Range("Q9").Value = 2.286358752355556
$ws.Range("R9").Value = 20.5772287712
$ws.Range("S9").Value = 0.0004216940590688622
$ws.Range("T9").Value = 0.0004216940590688622

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.08476666666667
$ws.Range("H10").Value = 33.2543
$ws.Range("I10").Value = 0.1808492350906109
$ws.Range("J10").Value = 0.1808492350906109
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.881585
$ws.Range("N10").Value = 5.644755
$ws.Range("O10").Value = 0.02127094711161878
$ws.Range("P10").Value = 0.02127094711161878
$ws.Range("Q10").Value = 20.8569306885
$ws.Range("R10").Value = 187.7123761965
$ws.Range("S10").Value = 0.003846834514789095
$ws.Range("T10").Value = 0.003846834514789095

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 11.08476666666667
$ws.Range("H11").Value = 33.2543
$ws.Range("I11").Value = 0.1808492350906109
$ws.Range("J11").Value = 0.1808492350906109
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2952356666666667
$ws.Range("N11").Value = 0.885707
$ws.Range("O11").Value = 0.003337580949640955
$ws.Range("P11").Value = 0.003337580949640955
$ws.Range("Q11").Value = 3.272618476677778
$ws.Range("R11").Value = 29.4535662901
$ws.Range("S11").Value = 0.0006035989617955615
$ws.Range("T11").Value = 0.0006035989617955615

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.140061
$ws.Range("H12").Value = 0.420183
$ws.Range("I12").Value = 0.002285111223152439
$ws.Range("J12").Value = 0.002285111223152439
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 83.40125033333334
$ws.Range("N12").Value = 250.203751
$ws.Range("O12").Value = 0.9428346765536562
$ws.Range("P12").Value = 0.9428346765536562
$ws.Range("Q12").Value = 11.681262522937
$ws.Range("R12").Value = 105.131362706433
$ws.Range("S12").Value = 0.00215448210097006
$ws.Range("T12").Value = 0.00215448210097006

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.140061
$ws.Range("H13").Value = 0.420183
$ws.Range("I13").Value = 0.002285111223152439
$ws.Range("J13").Value = 0.002285111223152439
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.673647
$ws.Range("N13").Value = 8.020941
$ws.Range("O13").Value = 0.03022505171551549
$ws.Range("P13").Value = 0.03022505171551549
$ws.Range("Q13").Value = 0.374473672467
$ws.Range("R13").Value = 3.370263052203
$ws.Range("S13").Value = 0.00006906760489548732
$ws.Range("T13").Value = 0.00006906760489548732

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.140061
$ws.Range("H14").Value = 0.420183
$ws.Range("I14").Value = 0.002285111223152439
$ws.Range("J14").Value = 0.002285111223152439
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2062613333333333
$ws.Range("N14").Value = 0.618784
$ws.Range("O14").Value = 0.002331743669568637
$ws.Range("P14").Value = 0.002331743669568637
$ws.Range("Q14").Value = 0.028889168608
$ws.Range("R14").Value = 0.260002517472
$ws.Range("S14").Value = 0.000005328293628845945
$ws.Range("T14").Value = 0.000005328293628845945

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.140061
$ws.Range("H15").Value = 0.420183
$ws.Range("I15").Value = 0.002285111223152439
$ws.Range("J15").Value = 0.002285111223152439
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.881585
$ws.Range("N15").Value = 5.644755
$ws.Range("O15").Value = 0.02127094711161878
$ws.Range("P15").Value = 0.02127094711161878
$ws.Range("Q15").Value = 0.263536676685
$ws.Range("R15").Value = 2.371830090165
$ws.Range("S15").Value = 0.00004860647997184202
$ws.Range("T15").Value = 0.00004860647997184202

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.140061
$ws.Range("H16").Value = 0.420183
$ws.Range("I16").Value = 0.002285111223152439
$ws.Range("J16").Value = 0.002285111223152439
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2952356666666667
$ws.Range("N16").Value = 0.885707
$ws.Range("O16").Value = 0.003337580949640955
$ws.Range("P16").Value = 0.003337580949640955
$ws.Range("Q16").Value = 0.041351002709
$ws.Range("R16").Value = 0.372159024381
$ws.Range("S16").Value = 0.000007626743686204322
$ws.Range("T16").Value = 0.000007626743686204323

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 7.750570000000001
$ws.Range("H17").Value = 23.25171
$ws.Range("I17").Value = 0.1264514353947823
$ws.Range("J17").Value = 0.1264514353947823
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 83.40125033333334
$ws.Range("N17").Value = 250.203751
$ws.Range("O17").Value = 0.9428346765536562
$ws.Range("P17").Value = 0.9428346765536562
$ws.Range("Q17").Value = 646.4072287960234
$ws.Range("R17").Value = 5817.665059164211
$ws.Range("S17").Value = 0.1192227981901851
$ws.Range("T17").Value = 0.1192227981901851

$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 7.750570000000001
$ws.Range("H18").Value = 23.25171
$ws.Range("I18").Value = 0.1264514353947823
$ws.Range("J18").Value = 0.1264514353947823
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 2.673647
$ws.Range("N18").Value = 8.020941
$ws.Range("O18").Value = 0.03022505171551549
$ws.Range("P18").Value = 0.03022505171551549
$ws.Range("Q18").Value = 20.72228822879
$ws.Range("R18").Value = 186.50059405911
$ws.Range("S18").Value = 0.00382200117430846
$ws.Range("T18").Value = 0.00382200117430846

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 7.750570000000001
$ws.Range("H19").Value = 23.25171
$ws.Range("I19").Value = 0.1264514353947823
$ws.Range("J19").Value = 0.1264514353947823
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.2062613333333333
$ws.Range("N19").Value = 0.618784
$ws.Range("O19").Value = 0.002331743669568637
$ws.Range("P19").Value = 0.002331743669568637
$ws.Range("Q19").Value = 1.598642902293334
$ws.Range("R19").Value = 14.38778612064
$ws.Range("S19").Value = 0.0002948523339896511
$ws.Range("T19").Value = 0.0002948523339896511

$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 7.750570000000001
$ws.Range("H20").Value = 23.25171
$ws.Range("I20").Value = 0.1264514353947823
$ws.Range("J20").Value = 0.1264514353947823
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 1.881585
$ws.Range("N20").Value = 5.644755
$ws.Range("O20").Value = 0.02127094711161878
$ws.Range("P20").Value = 0.02127094711161878
$ws.Range("Q20").Value = 14.58335625345
$ws.Range("R20").Value = 131.25020628105
$ws.Range("S20").Value = 0.002689741794470693
$ws.Range("T20").Value = 0.002689741794470693

$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 7.750570000000001
$ws.Range("H21").Value = 23.25171
$ws.Range("I21").Value = 0.1264514353947823
$ws.Range("J21").Value = 0.1264514353947823
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2952356666666667
$ws.Range("N21").Value = 0.885707
$ws.Range("O21").Value = 0.003337580949640955
$ws.Range("P21").Value = 0.003337580949640955
$ws.Range("Q21").Value = 2.288244700996667
$ws.Range("R21").Value = 20.59420230897
$ws.Range("S21").Value = 0.0004220419018283793
$ws.Range("T21").Value = 0.0004220419018283794
